# Auto-generated edit applying the scheduled-runner price/profit updates
# across the Tonberry_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

function Set-CellValue($ws, $row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = $value
}

function Clear-CellValue($ws, $row, $col) {
    $ws.Cells.Item($row, $col).ClearContents()
}

$ws = $wb.Worksheets.Item("ALC")
Set-CellValue $ws 32 8 1101
Set-CellValue $ws 32 9 1000
Set-CellValue $ws 32 10 1151.5
Set-CellValue $ws 32 11 1000
Set-CellValue $ws 32 12 1151.5
Set-CellValue $ws 32 13 -674
Set-CellValue $ws 32 14 -1803.5

$ws = $wb.Worksheets.Item("ALC")
Set-CellValue $ws 41 8 713.3333
Set-CellValue $ws 41 10 769.3077
Set-CellValue $ws 41 12 769.3077
Set-CellValue $ws 41 14 -1649.3077

$ws = $wb.Worksheets.Item("ALC")
Set-CellValue $ws 70 8 1521
Set-CellValue $ws 70 10 1521
Set-CellValue $ws 70 12 4563
Set-CellValue $ws 70 14 -5103

$ws = $wb.Worksheets.Item("ALC")
Set-CellValue $ws 73 8 1521
Set-CellValue $ws 73 10 1521
Set-CellValue $ws 73 12 4563
Set-CellValue $ws 73 14 -6435

$ws = $wb.Worksheets.Item("ALC")
Set-CellValue $ws 86 8 1232133
Set-CellValue $ws 86 10 0
Set-CellValue $ws 86 12 0
Clear-CellValue $ws 86 14

$ws = $wb.Worksheets.Item("ALC")
Set-CellValue $ws 89 8 1232133
Set-CellValue $ws 89 10 0
Set-CellValue $ws 89 12 0
Clear-CellValue $ws 89 14

$ws = $wb.Worksheets.Item("ALC")
Set-CellValue $ws 106 8 1406.091
Set-CellValue $ws 106 9 1406.091
Set-CellValue $ws 106 11 1406.091
Set-CellValue $ws 106 13 -775.0909999999999

$ws = $wb.Worksheets.Item("ALC")
Set-CellValue $ws 129 8 875.3953
Set-CellValue $ws 129 10 891.10254
Set-CellValue $ws 129 12 2673.30762
Set-CellValue $ws 129 14 -12673.30762

$ws = $wb.Worksheets.Item("ALC")
Set-CellValue $ws 137 8 2252.1365
Set-CellValue $ws 137 9 1338.3
Set-CellValue $ws 137 10 3013.6667
Set-CellValue $ws 137 11 4014.9
Set-CellValue $ws 137 12 9041.000100000001
Set-CellValue $ws 137 13 -1464.9
Set-CellValue $ws 137 14 -14141.0001

$ws = $wb.Worksheets.Item("ARM")
Set-CellValue $ws 32 8 3470.2
Set-CellValue $ws 32 9 3129.2046
Set-CellValue $ws 32 11 3129.2046
Set-CellValue $ws 32 13 -2842.2046

$ws = $wb.Worksheets.Item("ARM")
Set-CellValue $ws 61 8 4519.7
Set-CellValue $ws 61 9 1219.8
Set-CellValue $ws 61 10 7819.6
Set-CellValue $ws 61 11 1219.8
Set-CellValue $ws 61 12 7819.6
Set-CellValue $ws 61 13 -1007.8
Set-CellValue $ws 61 14 -8243.6

$ws = $wb.Worksheets.Item("ARM")
Set-CellValue $ws 123 8 63998
Set-CellValue $ws 123 10 63998
Set-CellValue $ws 123 12 63998
Set-CellValue $ws 123 14 -73798

$ws = $wb.Worksheets.Item("ARM")
Set-CellValue $ws 132 8 1654
Set-CellValue $ws 132 9 1026.5312
Set-CellValue $ws 132 11 3079.5936
Set-CellValue $ws 132 13 -549.5935999999997

$ws = $wb.Worksheets.Item("ARM")
Set-CellValue $ws 136 8 4519.7
Set-CellValue $ws 136 9 1219.8
Set-CellValue $ws 136 10 7819.6
Set-CellValue $ws 136 11 3659.4
Set-CellValue $ws 136 12 23458.8
Set-CellValue $ws 136 13 -1109.4
Set-CellValue $ws 136 14 -28558.8

$ws = $wb.Worksheets.Item("BSM")
Set-CellValue $ws 20 8 2055.4
Set-CellValue $ws 20 9 1982.7307
Set-CellValue $ws 20 10 2527.75
Set-CellValue $ws 20 11 1982.7307
Set-CellValue $ws 20 12 2527.75
Set-CellValue $ws 20 13 -1735.7307
Set-CellValue $ws 20 14 -3021.75

$ws = $wb.Worksheets.Item("BSM")
Set-CellValue $ws 94 8 786.1539
Set-CellValue $ws 94 9 620
Set-CellValue $ws 94 11 620
Set-CellValue $ws 94 13 -169

$ws = $wb.Worksheets.Item("BSM")
Set-CellValue $ws 105 8 2548.4
Set-CellValue $ws 105 9 2435.5625
Set-CellValue $ws 105 10 2999.75
Set-CellValue $ws 105 11 2435.5625
Set-CellValue $ws 105 12 2999.75
Set-CellValue $ws 105 13 -688.5625
Set-CellValue $ws 105 14 -6493.75

$ws = $wb.Worksheets.Item("BSM")
Set-CellValue $ws 134 8 5940.2856
Set-CellValue $ws 134 9 6571.857
Set-CellValue $ws 134 11 19715.571
Set-CellValue $ws 134 13 -17180.571

$ws = $wb.Worksheets.Item("CRP")
Set-CellValue $ws 31 8 2481.7334
Set-CellValue $ws 31 9 2448.111
Set-CellValue $ws 31 11 2448.111
Set-CellValue $ws 31 13 -2153.111

$ws = $wb.Worksheets.Item("CRP")
Set-CellValue $ws 33 8 5890
Set-CellValue $ws 33 9 3900
Set-CellValue $ws 33 10 9870
Set-CellValue $ws 33 11 3900
Set-CellValue $ws 33 12 9870
Set-CellValue $ws 33 13 -3521
Set-CellValue $ws 33 14 -10628

$ws = $wb.Worksheets.Item("CRP")
Set-CellValue $ws 34 8 2481.7334
Set-CellValue $ws 34 9 2448.111
Set-CellValue $ws 34 11 2448.111
Set-CellValue $ws 34 13 -2246.111

$ws = $wb.Worksheets.Item("CRP")
Set-CellValue $ws 122 8 3922.8
Set-CellValue $ws 122 9 1771.5714
Set-CellValue $ws 122 11 5314.7142
Set-CellValue $ws 122 13 -2864.7142

$ws = $wb.Worksheets.Item("CRP")
Set-CellValue $ws 141 8 71998
Set-CellValue $ws 141 10 71998
Set-CellValue $ws 141 12 71998
Set-CellValue $ws 141 14 -82358

$ws = $wb.Worksheets.Item("CUL")
Set-CellValue $ws 20 8 850
Set-CellValue $ws 20 9 700
Set-CellValue $ws 20 11 2100
Set-CellValue $ws 20 13 -1873

$ws = $wb.Worksheets.Item("CUL")
Set-CellValue $ws 29 8 105.57143
Set-CellValue $ws 29 9 72.5
Set-CellValue $ws 29 10 118.8
Set-CellValue $ws 29 11 217.5
Set-CellValue $ws 29 12 356.4
Set-CellValue $ws 29 13 59.5
Set-CellValue $ws 29 14 -910.4

$ws = $wb.Worksheets.Item("CUL")
Set-CellValue $ws 108 8 0
Set-CellValue $ws 108 9 0
Set-CellValue $ws 108 10 0
Set-CellValue $ws 108 11 0
Set-CellValue $ws 108 12 0
Clear-CellValue $ws 108 13
Clear-CellValue $ws 108 14

$ws = $wb.Worksheets.Item("CUL")
Set-CellValue $ws 130 8 250000660
Set-CellValue $ws 130 9 250000660
Set-CellValue $ws 130 10 0
Set-CellValue $ws 130 11 750001980
Set-CellValue $ws 130 12 0
Set-CellValue $ws 130 13 -749996960
Clear-CellValue $ws 130 14

$ws = $wb.Worksheets.Item("GSM")
Set-CellValue $ws 122 8 1506.4667
Set-CellValue $ws 122 9 975.125
Set-CellValue $ws 122 10 2113.7144
Set-CellValue $ws 122 11 2925.375
Set-CellValue $ws 122 12 6341.1432
Set-CellValue $ws 122 13 -475.375
Set-CellValue $ws 122 14 -11241.1432

$ws = $wb.Worksheets.Item("GSM")
Set-CellValue $ws 132 8 1328795.1
Set-CellValue $ws 132 9 1924875.5
Set-CellValue $ws 132 10 4172.1113
Set-CellValue $ws 132 11 5774626.5
Set-CellValue $ws 132 12 12516.3339
Set-CellValue $ws 132 13 -5772096.5
Set-CellValue $ws 132 14 -17576.3339

$ws = $wb.Worksheets.Item("LTW")
Set-CellValue $ws 61 8 1657.625
Set-CellValue $ws 61 9 1607.9375
Set-CellValue $ws 61 10 1757
Set-CellValue $ws 61 11 1607.9375
Set-CellValue $ws 61 12 1757
Set-CellValue $ws 61 13 -1405.9375
Set-CellValue $ws 61 14 -2161

$ws = $wb.Worksheets.Item("LTW")
Set-CellValue $ws 100 8 1213.2858
Set-CellValue $ws 100 9 1213.2858
Set-CellValue $ws 100 10 0
Set-CellValue $ws 100 11 1213.2858
Set-CellValue $ws 100 12 0
Set-CellValue $ws 100 13 -672.2858000000001
Clear-CellValue $ws 100 14

$ws = $wb.Worksheets.Item("LTW")
Set-CellValue $ws 113 8 1657.625
Set-CellValue $ws 113 9 1607.9375
Set-CellValue $ws 113 10 1757
Set-CellValue $ws 113 11 1607.9375
Set-CellValue $ws 113 12 1757
Set-CellValue $ws 113 13 562.0625
Set-CellValue $ws 113 14 -6097

$ws = $wb.Worksheets.Item("LTW")
Set-CellValue $ws 122 8 4309.4
Set-CellValue $ws 122 9 2847.3333
Set-CellValue $ws 122 11 8541.999899999999
Set-CellValue $ws 122 13 -6091.999899999999

$ws = $wb.Worksheets.Item("WVR")
Set-CellValue $ws 14 8 1430.44
Set-CellValue $ws 14 9 998
Set-CellValue $ws 14 10 1468.0435
Set-CellValue $ws 14 11 998
Set-CellValue $ws 14 12 1468.0435
Set-CellValue $ws 14 13 -830
Set-CellValue $ws 14 14 -1804.0435

$ws = $wb.Worksheets.Item("WVR")
Set-CellValue $ws 100 8 499.85715
Set-CellValue $ws 100 9 383.16666
Set-CellValue $ws 100 10 1200
Set-CellValue $ws 100 11 766.33332
Set-CellValue $ws 100 12 2400
Set-CellValue $ws 100 13 -225.33332
Set-CellValue $ws 100 14 -3482

$ws = $wb.Worksheets.Item("WVR")
Set-CellValue $ws 124 8 24996
Set-CellValue $ws 124 10 24996
Set-CellValue $ws 124 12 24996
Set-CellValue $ws 124 14 -34816
